$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows whose Target cluster (column D) is "ECs" (shared string index 20).
# These are the original rows 2, 5, 8 and 11; deleted from the bottom up so the
# remaining row numbers do not shift under us while we are iterating.
$ws.Rows.Item(11).EntireRow.Delete()
$ws.Rows.Item(8).EntireRow.Delete()
$ws.Rows.Item(5).EntireRow.Delete()
$ws.Rows.Item(2).EntireRow.Delete()

# Refresh the recalculated NATMI TPM statistics for the remaining sender/target rows.
$ws.Range("G2").Value = 6.221859333333334
$ws.Range("H2").Value = 18.665578
$ws.Range("I2").Value = 0.266168288812213
$ws.Range("J2").Value = 0.2661682888122131
$ws.Range("M2").Value = 0.8467083333333334
$ws.Range("N2").Value = 2.540125
$ws.Range("O2").Value = 0.1482255704769734
$ws.Range("P2").Value = 0.1482255704769734
$ws.Range("Q2").Value = 5.268100146361112
$ws.Range("R2").Value = 47.41290131725
$ws.Range("S2").Value = 0.0394529464520701
$ws.Range("T2").Value = 0.03945294645207011

$ws.Range("G3").Value = 6.221859333333334
$ws.Range("H3").Value = 18.665578
$ws.Range("I3").Value = 0.266168288812213
$ws.Range("J3").Value = 0.2661682888122131
$ws.Range("M3").Value = 4.865587666666666
$ws.Range("N3").Value = 14.596763
$ws.Range("O3").Value = 0.8517744295230265
$ws.Range("P3").Value = 0.8517744295230266
$ws.Range("Q3").Value = 30.27300203600155
$ws.Range("R3").Value = 272.457018324014
$ws.Range("S3").Value = 0.2267153423601429
$ws.Range("T3").Value = 0.226715342360143

$ws.Range("G4").Value = 6.924657666666666
$ws.Range("I4").Value = 0.2962336791949928
$ws.Range("J4").Value = 0.2962336791949928
$ws.Range("M4").Value = 0.8467083333333334
$ws.Range("N4").Value = 2.540125
$ws.Range("O4").Value = 0.1482255704769734
$ws.Range("P4").Value = 0.1482255704769734
$ws.Range("Q4").Value = 5.863165351847222
$ws.Range("R4").Value = 52.768488166625
$ws.Range("S4").Value = 0.04390940609317053
$ws.Range("T4").Value = 0.04390940609317055

$ws.Range("G5").Value = 6.924657666666666
$ws.Range("I5").Value = 0.2962336791949928
$ws.Range("J5").Value = 0.2962336791949928
$ws.Range("M5").Value = 4.865587666666666
$ws.Range("N5").Value = 14.596763
$ws.Range("O5").Value = 0.8517744295230265
$ws.Range("P5").Value = 0.8517744295230266
$ws.Range("Q5").Value = 33.69252893882211
$ws.Range("R5").Value = 303.2327604493989
$ws.Range("S5").Value = 0.2523242731018222
$ws.Range("T5").Value = 0.2523242731018223

$ws.Range("G6").Value = 2.674426
$ws.Range("H6").Value = 8.023277999999999
$ws.Range("I6").Value = 0.1144107177353241
$ws.Range("J6").Value = 0.1144107177353241
$ws.Range("M6").Value = 0.8467083333333334
$ws.Range("N6").Value = 2.540125
$ws.Range("O6").Value = 0.1482255704769734
$ws.Range("P6").Value = 0.1482255704769734
$ws.Range("Q6").Value = 2.264458781083333
$ws.Range("R6").Value = 20.38012902975
$ws.Range("S6").Value = 0.01695859390499839
$ws.Range("T6").Value = 0.0169585939049984

$ws.Range("G7").Value = 2.674426
$ws.Range("H7").Value = 8.023277999999999
$ws.Range("I7").Value = 0.1144107177353241
$ws.Range("J7").Value = 0.1144107177353241
$ws.Range("M7").Value = 4.865587666666666
$ws.Range("N7").Value = 14.596763
$ws.Range("O7").Value = 0.8517744295230265
$ws.Range("P7").Value = 0.8517744295230266
$ws.Range("Q7").Value = 13.01265416101267
$ws.Range("R7").Value = 117.113887449114
$ws.Range("S7").Value = 0.09745212383032568
$ws.Range("T7").Value = 0.09745212383032571

$ws.Range("G8").Value = 7.554716666666667
$ws.Range("H8").Value = 22.66415
$ws.Range("I8").Value = 0.32318731425747
$ws.Range("J8").Value = 0.32318731425747
$ws.Range("M8").Value = 0.8467083333333334
$ws.Range("N8").Value = 2.540125
$ws.Range("O8").Value = 0.1482255704769734
$ws.Range("P8").Value = 0.1482255704769734
$ws.Range("Q8").Value = 6.39664155763889
$ws.Range("R8").Value = 57.56977401875
$ws.Range("S8").Value = 0.04790462402673437
$ws.Range("T8").Value = 0.04790462402673439

$ws.Range("G9").Value = 7.554716666666667
$ws.Range("H9").Value = 22.66415
$ws.Range("I9").Value = 0.32318731425747
$ws.Range("J9").Value = 0.32318731425747
$ws.Range("M9").Value = 4.865587666666666
$ws.Range("N9").Value = 14.596763
$ws.Range("O9").Value = 0.8517744295230265
$ws.Range("P9").Value = 0.8517744295230266
$ws.Range("Q9").Value = 36.75813623849444
$ws.Range("R9").Value = 330.82322614645
$ws.Range("S9").Value = 0.2752826902307356
$ws.Range("T9").Value = 0.2752826902307357
